# Insert a new data row above row 25 (shifting existing rows 25-35 down to 26-36)
# and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("25:25").Insert()

$ws.Range("A25").Value = 10
$ws.Range("B25").Value = "Vega Modelo de Temuco"
$ws.Range("C25").Value = "La Araucanía"
$ws.Range("D25").Value = 44529
$ws.Range("E25").Value = 9
$ws.Range("F25").Value = 100114002
$ws.Range("G25").Value = "Camote"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 15
$ws.Range("K25").Value = 20000
$ws.Range("L25").Value = 20000
$ws.Range("M25").Value = 20000
$ws.Range("N25").Value = "`$/malla 20 kilos"
$ws.Range("O25").Value = "Perú"
$ws.Range("P25").Value = 1000
$ws.Range("Q25").Value = 20
$ws.Range("R25").Value = "Hortaliza"
